$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.693.20'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '3.042.76'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'379.76"
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = "'103.10"
$ws.Range("E6").Value = '  +1.67%  '
$ws.Range("E7").Value = '  +0.87%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("D10").Value = "'36.97"
$ws.Range("E10").Value = '  +2.05%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = "'0.0862"
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("D13").Value = '3.521.30'
$ws.Range("E13").Value = '  +3.15%  '
$ws.Range("D14").Value = "'18.60"
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = '3.034.29'
$ws.Range("E16").Value = '  +3.21%  '
$ws.Range("D17").Value = "'0.980"
$ws.Range("D18").Value = "'10.54"
$ws.Range("E18").Value = '  -11.37%  '
$ws.Range("D19").Value = '51.696.23'
$ws.Range("E19").Value = '  +1.38%  '
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("E21").Value = '  +0.61%  '
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").Value = "'69.99"
$ws.Range("D24").Value = "'269.11"
$ws.Range("E24").Value = '  +0.85%  '
$ws.Range("D25").Value = "'3.16"
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").Value = "'7.55"
$ws.Range("E27").Value = '  +6.57%  '
$ws.Range("E28").Value = '  +6.26%  '
$ws.Range("D29").Value = "'26.37"
$ws.Range("E29").Value = '  +2.91%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = "'10.33"
$ws.Range("E32").Value = '  +2.25%  '
$ws.Range("E33").Value = '  +1.91%  '
$ws.Range("D34").Value = "'50.48"
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = "'0.0451"
$ws.Range("E36").Value = '  +4.61%  '
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = '  +7.43%  '
$ws.Range("D39").Value = "'0.289"
$ws.Range("E39").Value = '  +10.84%  '
$ws.Range("D40").Value = "'17.11"
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("E41").Value = '  +3.44%  '
$ws.Range("D42").Value = "'2.59"
$ws.Range("E42").Value = '  +2.41%  '
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").Value = "'127.53"
$ws.Range("E44").Value = '  +8.30%  '
$ws.Range("D45").Value = "'3.73"
$ws.Range("E45").Value = '  +5.94%  '
$ws.Range("D46").Value = "'21.94"
$ws.Range("E46").Value = '  +2.70%  '
$ws.Range("D47").Value = "'2.14"
$ws.Range("E47").Value = '  +6.15%  '
$ws.Range("E48").Value = '  +3.02%  '
$ws.Range("D49").Value = '2.037.64'
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("D50").Value = '3.337.86'
$ws.Range("E50").Value = '  +2.96%  '
$ws.Range("D51").Value = "'0.0320"
$ws.Range("E51").Value = '  +2.26%  '
